$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple G-column value corrections (rows scattered through the sheet) ---
$ws.Cells.Item(56, 7).Value = "0,00"
$ws.Cells.Item(57, 7).Value = "0,00"
$ws.Cells.Item(58, 7).Value = "0,00"
$ws.Cells.Item(113, 7).Value = "0,00"
$ws.Cells.Item(115, 7).Value = "0,00"
$ws.Cells.Item(205, 7).Value = "1,00"
$ws.Cells.Item(225, 7).Value = "0,00"
$ws.Cells.Item(238, 7).Value = "1,00"
$ws.Cells.Item(240, 7).Value = "5,00"
$ws.Cells.Item(241, 7).Value = "0,00"
$ws.Cells.Item(242, 7).Value = "0,00"
$ws.Cells.Item(336, 7).Value = "0,00"
$ws.Cells.Item(337, 7).Value = "0,00"
$ws.Cells.Item(338, 7).Value = "0,00"
$ws.Cells.Item(339, 7).Value = "0,00"
$ws.Cells.Item(423, 7).Value = "0,00"
$ws.Cells.Item(458, 7).Value = "0,00"
$ws.Cells.Item(462, 7).Value = "0,00"
$ws.Cells.Item(463, 7).Value = "0,00"
$ws.Cells.Item(464, 7).Value = "0,00"
$ws.Cells.Item(469, 7).Value = "0,00"
$ws.Cells.Item(471, 7).Value = "0,00"
$ws.Cells.Item(517, 7).Value = "25,00"
$ws.Cells.Item(686, 7).Value = "0,00"
$ws.Cells.Item(939, 7).Value = "200,00"
$ws.Cells.Item(947, 7).Value = "0,00"
$ws.Cells.Item(985, 7).Value = "0,00"
$ws.Cells.Item(986, 7).Value = "0,00"
$ws.Cells.Item(987, 7).Value = "0,00"
$ws.Cells.Item(988, 7).Value = "60,85"
$ws.Cells.Item(989, 7).NumberFormat = "@"
$ws.Cells.Item(989, 7).Value = "337,956"
$ws.Cells.Item(989, 7).ClearFormats()
$ws.Cells.Item(990, 7).NumberFormat = "@"
$ws.Cells.Item(990, 7).Value = "555,298"
$ws.Cells.Item(990, 7).ClearFormats()
$ws.Cells.Item(1017, 7).Value = "322,00"
$ws.Cells.Item(1018, 7).Value = "0,00"
$ws.Cells.Item(1019, 7).Value = "0,00"
$ws.Cells.Item(1027, 7).Value = "2,76"
$ws.Cells.Item(1028, 7).Value = "0,00"
$ws.Cells.Item(1029, 7).Value = "0,00"
$ws.Cells.Item(1030, 7).Value = "275,21"
$ws.Cells.Item(1031, 7).Value = "4212,77"
$ws.Cells.Item(1032, 7).Value = "582,43"
$ws.Cells.Item(1034, 7).Value = "0,5"
$ws.Cells.Item(1055, 7).Value = "2,76"
$ws.Cells.Item(1056, 7).Value = "0,00"
$ws.Cells.Item(1057, 7).Value = "0,00"
$ws.Cells.Item(1058, 7).Value = "275,21"
$ws.Cells.Item(1059, 7).Value = "4212,77"
$ws.Cells.Item(1060, 7).Value = "582,43"
$ws.Cells.Item(1062, 7).Value = "0,5"
$ws.Cells.Item(1117, 7).Value = "1,04"
$ws.Cells.Item(1118, 7).Value = "0,00"
$ws.Cells.Item(1119, 7).Value = "0,00"
$ws.Cells.Item(1120, 7).Value = "5,84"
$ws.Cells.Item(1121, 7).Value = "772,64"
$ws.Cells.Item(1122, 7).Value = "114,2"
$ws.Cells.Item(1147, 7).Value = "0,00"
$ws.Cells.Item(1150, 7).Value = "0,00"
$ws.Cells.Item(1151, 7).Value = "10,00"
$ws.Cells.Item(1152, 7).Value = "1,00"
$ws.Cells.Item(1227, 7).Value = "164,01"
$ws.Cells.Item(1229, 7).Value = "164,01"
$ws.Cells.Item(1232, 7).Value = "164,01"
$ws.Cells.Item(1234, 7).Value = "164,01"
$ws.Cells.Item(1238, 7).Value = "1,00"
$ws.Cells.Item(1239, 7).Value = "14,00"
$ws.Cells.Item(1240, 7).Value = "1,00"
$ws.Cells.Item(1253, 7).Value = "181,44"
$ws.Cells.Item(1254, 7).Value = "15,6"
$ws.Cells.Item(1257, 7).NumberFormat = "@"
$ws.Cells.Item(1257, 7).Value = "983,325"
$ws.Cells.Item(1257, 7).ClearFormats()
$ws.Cells.Item(1499, 7).Value = "0,00"
$ws.Cells.Item(1594, 7).Value = "180,33"

# --- Rebuild rows 1600-1635 (item catalog rows shifted / replaced / appended) ---
$ws.Cells.Item(1600, 1).Value = "-"
$ws.Cells.Item(1600, 2).Value = "-"
$ws.Cells.Item(1600, 3).Value = "-"
$ws.Cells.Item(1600, 4).Value = "-"
$ws.Cells.Item(1600, 5).Value = "SUMINISTRO DE ELBOW 45° LR, SCH 40, BW, ASTM A-234 GR WPB, ASME B16.9 6"""
$ws.Cells.Item(1600, 6).Value = "UN"
$ws.Cells.Item(1600, 7).Value = "3,00"
$ws.Cells.Item(1600, 8).Value = "NUEVO"

$ws.Cells.Item(1601, 1).Value = "-"
$ws.Cells.Item(1601, 2).Value = "-"
$ws.Cells.Item(1601, 3).Value = "-"
$ws.Cells.Item(1601, 4).Value = "-"
$ws.Cells.Item(1601, 5).Value = "SUMINISTRO DE 90 DEGREE ELBOW LR ASME B16.9 BW WROUGHT S ASTM A234 GR. WPB, 4"", S-30"
$ws.Cells.Item(1601, 6).Value = "UN"
$ws.Cells.Item(1601, 7).Value = "7,00"
$ws.Cells.Item(1601, 8).Value = "NUEVO"

$ws.Cells.Item(1602, 1).Value = "-"
$ws.Cells.Item(1602, 2).Value = "-"
$ws.Cells.Item(1602, 3).Value = "-"
$ws.Cells.Item(1602, 4).Value = "-"
$ws.Cells.Item(1602, 5).Value = "SUMINISTRO DE 90 DEGREE ELBOW LR ASME B16.9 BW WROUGHT S ASTM A234 GR. WPB, 6"", S-40"
$ws.Cells.Item(1602, 6).Value = "UN"
$ws.Cells.Item(1602, 7).Value = "19,00"
$ws.Cells.Item(1602, 8).Value = "NUEVO"

$ws.Cells.Item(1603, 1).Value = "-"
$ws.Cells.Item(1603, 2).Value = "-"
$ws.Cells.Item(1603, 3).Value = "-"
$ws.Cells.Item(1603, 4).Value = "-"
$ws.Cells.Item(1603, 5).Value = "SUMINISTRO DE PIPE ASME B36.10 BE SEAMLESS ASTM A53 GR.  B, 2-1/2 "", S-30"
$ws.Cells.Item(1603, 6).Value = "M"
$ws.Cells.Item(1603, 7).Value = "1,00"
$ws.Cells.Item(1603, 8).Value = "NUEVO"

$ws.Cells.Item(1604, 1).Value = "-"
$ws.Cells.Item(1604, 2).Value = "-"
$ws.Cells.Item(1604, 3).Value = "-"
$ws.Cells.Item(1604, 4).Value = "-"
$ws.Cells.Item(1604, 5).Value = "SUMINISTRO DE PIPE ASME B36.10 BE SEAMLESS ASTM A53 GR.  B, 4 "", S-30"
$ws.Cells.Item(1604, 6).Value = "M"
$ws.Cells.Item(1604, 7).Value = "17,00"
$ws.Cells.Item(1604, 8).Value = "NUEVO"

$ws.Cells.Item(1605, 1).Value = "-"
$ws.Cells.Item(1605, 2).Value = "-"
$ws.Cells.Item(1605, 3).Value = "-"
$ws.Cells.Item(1605, 4).Value = "-"
$ws.Cells.Item(1605, 5).Value = "SUMINISTRO DE PIPE ASME B36.10 BE SEAMLESS ASTM A53 GR.  B, 6 "", S-40"
$ws.Cells.Item(1605, 6).Value = "M"
$ws.Cells.Item(1605, 7).Value = "127,00"
$ws.Cells.Item(1605, 8).Value = "NUEVO"

$ws.Cells.Item(1606, 1).Value = "-"
$ws.Cells.Item(1606, 2).Value = "-"
$ws.Cells.Item(1606, 3).Value = "-"
$ws.Cells.Item(1606, 4).Value = "-"
$ws.Cells.Item(1606, 5).Value = "SUMINISTRO DE PIPE ASME B36.10 BE SEAMLESS ASTM A53 GR.  B, 8 "", S-20"
$ws.Cells.Item(1606, 6).Value = "M"
$ws.Cells.Item(1606, 7).Value = "13,00"
$ws.Cells.Item(1606, 8).Value = "NUEVO"

$ws.Cells.Item(1607, 1).Value = "-"
$ws.Cells.Item(1607, 2).Value = "-"
$ws.Cells.Item(1607, 3).Value = "-"
$ws.Cells.Item(1607, 4).Value = "-"
$ws.Cells.Item(1607, 5).Value = "SUMINISTRO BLIND FLANGE, CLASS 150, FF, ASTM A105, ASME B16.5. DIAMETRO DE 6"""
$ws.Cells.Item(1607, 6).Value = "UN"
$ws.Cells.Item(1607, 7).Value = "2,00"
$ws.Cells.Item(1607, 8).Value = "NUEVO"

$ws.Cells.Item(1608, 1).Value = "-"
$ws.Cells.Item(1608, 2).Value = "-"
$ws.Cells.Item(1608, 3).Value = "-"
$ws.Cells.Item(1608, 4).Value = "-"
$ws.Cells.Item(1608, 5).Value = "SUMINISTRO DE STUD BOLTS, A-193 GR B7 W/(2) HVY NUTS HEX, ASTM A-194 GR 2H, ASME B18.2.1, B18.2.2 Ø3/4""X120MMLG"
$ws.Cells.Item(1608, 6).Value = "UN"
$ws.Cells.Item(1608, 7).Value = "24,00"
$ws.Cells.Item(1608, 8).Value = "NUEVO"

$ws.Cells.Item(1609, 1).Value = "-"
$ws.Cells.Item(1609, 2).Value = "-"
$ws.Cells.Item(1609, 3).Value = "-"
$ws.Cells.Item(1609, 4).Value = "-"
$ws.Cells.Item(1609, 5).Value = "SUMINISTRO DE STUD BOLTS, A-193 GR B7 W/(2) HVY NUTS HEX, ASTM A-194 GR 2H, ASME B18.2.1, B18.2.2 Ø5/8""X100MMLG"
$ws.Cells.Item(1609, 6).Value = "UN"
$ws.Cells.Item(1609, 7).Value = "24,00"
$ws.Cells.Item(1609, 8).Value = "NUEVO"

$ws.Cells.Item(1610, 1).Value = "-"
$ws.Cells.Item(1610, 2).Value = "-"
$ws.Cells.Item(1610, 3).Value = "-"
$ws.Cells.Item(1610, 4).Value = "-"
$ws.Cells.Item(1610, 5).Value = "SUMINISTRO DE SWING CHECK VALVE FL, MSS SP-136, A536 Gr 65-45-12, CL 125, INST HORIZ/VERT, FF, B16.1, BOLTED COVER, NON METALLIC FLAT GASKET EPDM, SEATS BRONZE, DISC DUCTILE IRON EPDM ENCAPSULATED, SS304 PIN,UL LISTED/FM APPROVED, 6"""
$ws.Cells.Item(1610, 6).Value = "UN"
$ws.Cells.Item(1610, 7).Value = "1,00"
$ws.Cells.Item(1610, 8).Value = "NUEVO"

$ws.Cells.Item(1611, 1).Value = "-"
$ws.Cells.Item(1611, 2).Value = "-"
$ws.Cells.Item(1611, 3).Value = "-"
$ws.Cells.Item(1611, 4).Value = "-"
$ws.Cells.Item(1611, 5).Value = "SUMINISTRO DE TEE, SCH 20, BW, ASTM A-234 GR WPB, ASME B16.9 12"""
$ws.Cells.Item(1611, 6).Value = "UN"
$ws.Cells.Item(1611, 7).Value = "2,00"
$ws.Cells.Item(1611, 8).Value = "NUEVO"

$ws.Cells.Item(1612, 1).Value = "-"
$ws.Cells.Item(1612, 2).Value = "-"
$ws.Cells.Item(1612, 3).Value = "-"
$ws.Cells.Item(1612, 4).Value = "-"
$ws.Cells.Item(1612, 5).Value = "SUMINISTRO DE CONC. REDUCER, SCH 40 X SCH 30, BW, ASTM A-234 GR WPB, ASME B16.9 6""X4"""
$ws.Cells.Item(1612, 6).Value = "UN"
$ws.Cells.Item(1612, 7).Value = "1,00"
$ws.Cells.Item(1612, 8).Value = "NUEVO"

$ws.Cells.Item(1613, 1).Value = "-"
$ws.Cells.Item(1613, 2).Value = "-"
$ws.Cells.Item(1613, 3).Value = "-"
$ws.Cells.Item(1613, 4).Value = "-"
$ws.Cells.Item(1613, 5).Value = "SUMINISTRO DE EXCEN. REDUCER, SCH 20 X SCH 40, BW, ASTM A-234 GR WPB, ASME B16.9 8""X6"""
$ws.Cells.Item(1613, 6).Value = "UN"
$ws.Cells.Item(1613, 7).Value = "1,00"
$ws.Cells.Item(1613, 8).Value = "NUEVO"

$ws.Cells.Item(1614, 1).Value = "-"
$ws.Cells.Item(1614, 2).Value = "-"
$ws.Cells.Item(1614, 3).Value = "-"
$ws.Cells.Item(1614, 4).Value = "-"
$ws.Cells.Item(1614, 5).Value = "SUMINISTRO DE NONMETALLIC FLAT GASKET, ARAMID FIBER REINFORCED NBR, ASME B16.21, THK 1/8"", CLASS 150, FF, 4"""
$ws.Cells.Item(1614, 6).Value = "UN"
$ws.Cells.Item(1614, 7).Value = "4,00"
$ws.Cells.Item(1614, 8).Value = "NUEVO"

$ws.Cells.Item(1615, 1).Value = "-"
$ws.Cells.Item(1615, 2).Value = "-"
$ws.Cells.Item(1615, 3).Value = "-"
$ws.Cells.Item(1615, 4).Value = "-"
$ws.Cells.Item(1615, 5).Value = "SUMINISTRO DE NONMETALLIC FLAT GASKET, ARAMID FIBER REINFORCED NBR, ASME B16.21, THK 1/8"", CLASS 150, FF, 6"""
$ws.Cells.Item(1615, 6).Value = "UN"
$ws.Cells.Item(1615, 7).Value = "24,00"
$ws.Cells.Item(1615, 8).Value = "NUEVO"

$ws.Cells.Item(1616, 1).Value = "-"
$ws.Cells.Item(1616, 2).Value = "-"
$ws.Cells.Item(1616, 3).Value = "-"
$ws.Cells.Item(1616, 4).Value = "-"
$ws.Cells.Item(1616, 5).Value = "SUMINISTRO DE NONMETALLIC FLAT GASKET, ARAMID FIBER REINFORCED NBR, ASME B16.21, THK 1/8"", CLASS 150, FF, 8"""
$ws.Cells.Item(1616, 6).Value = "UN"
$ws.Cells.Item(1616, 7).Value = "4,00"
$ws.Cells.Item(1616, 8).Value = "NUEVO"

$ws.Cells.Item(1617, 1).Value = "-"
$ws.Cells.Item(1617, 2).Value = "-"
$ws.Cells.Item(1617, 3).Value = "-"
$ws.Cells.Item(1617, 4).Value = "-"
$ws.Cells.Item(1617, 5).Value = "SUMINISTRO DE GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED, 6"""
$ws.Cells.Item(1617, 6).Value = "UN"
$ws.Cells.Item(1617, 7).Value = "9,00"
$ws.Cells.Item(1617, 8).Value = "NUEVO"

$ws.Cells.Item(1618, 1).Value = "-"
$ws.Cells.Item(1618, 2).Value = "-"
$ws.Cells.Item(1618, 3).Value = "-"
$ws.Cells.Item(1618, 4).Value = "-"
$ws.Cells.Item(1618, 5).Value = "SUMINISTRO DE GATE VALVE FL, MSS SP-128 TYPE II, A536 Gr 65-45-12, CL 125, FF, B16.1, BB, NON METALLIC FLAT GASKET EPDM, PKG EPDM, SOLID WEDGE, WEDGE DUCTILE IRON EPDM ENCAPSULATED, STEM BRONZE, OS&Y/RSNRO, HO, UL LISTED/FM APPROVED, 8"""
$ws.Cells.Item(1618, 6).Value = "UN"
$ws.Cells.Item(1618, 7).Value = "1,00"
$ws.Cells.Item(1618, 8).Value = "NUEVO"

$ws.Cells.Item(1619, 1).Value = "-"
$ws.Cells.Item(1619, 2).Value = "-"
$ws.Cells.Item(1619, 3).Value = "-"
$ws.Cells.Item(1619, 4).Value = "-"
$ws.Cells.Item(1619, 5).Value = "SUMINISTRO DE GATE VALVE THD, MSS SP-80, B62 UNS C83600, CL 150, NPTF, B1.20.1, SCREWED BONNET, PKG LUBRICATED FIBER/GRAPH, BRONZE STEM, SEAT RINGS&DISC, S, SOLID WEDGE, STEM OS&Y/RSNRO, HO, 1/2"""
$ws.Cells.Item(1619, 6).Value = "UN"
$ws.Cells.Item(1619, 7).Value = "1,00"
$ws.Cells.Item(1619, 8).Value = "NUEVO"

$ws.Cells.Item(1620, 1).Value = "-"
$ws.Cells.Item(1620, 2).Value = "-"
$ws.Cells.Item(1620, 3).Value = "-"
$ws.Cells.Item(1620, 4).Value = "-"
$ws.Cells.Item(1620, 5).Value = "SUMINISTRO DE ANGLE HOSE VALVE THD, MSS SP-80, B62 UNS C83600, 300 PSI CWP, FNPT X NH W/CAP AND CHAIN, FNPT B1.20.1 AND NH NFPA 1963, SCREW-IN BONNET, RENEWABLE DISC, PKG NON ASBESTOS, DISC&STEM BRONZE, RISING STEM, HO, UL LISTED/FM APPROVED, 2-1/2"""
$ws.Cells.Item(1620, 6).Value = "UN"
$ws.Cells.Item(1620, 7).Value = "4,00"
$ws.Cells.Item(1620, 8).Value = "NUEVO"

$ws.Cells.Item(1621, 1).Value = "-"
$ws.Cells.Item(1621, 2).Value = "-"
$ws.Cells.Item(1621, 3).Value = "-"
$ws.Cells.Item(1621, 4).Value = "-"
$ws.Cells.Item(1621, 5).Value = "SUMINISTRO DE PIPExTE NIPPLE - LENGTH = 100 MM ASME B36.10 PExTE SEAMLESS ASTM A53 GR. B, 1/2 "", S-XXS"
$ws.Cells.Item(1621, 6).Value = "UN"
$ws.Cells.Item(1621, 7).Value = "1,00"
$ws.Cells.Item(1621, 8).Value = "NUEVO"

$ws.Cells.Item(1622, 1).Value = "-"
$ws.Cells.Item(1622, 2).Value = "-"
$ws.Cells.Item(1622, 3).Value = "-"
$ws.Cells.Item(1622, 4).Value = "-"
$ws.Cells.Item(1622, 5).Value = "SUMINISTRO DE PIPExTE NIPPLE - LENGTH = 100 MM ASME B36.10 PExTE SEAMLESS ASTM A53 GR. B, 1½ "", S-160"
$ws.Cells.Item(1622, 6).Value = "UN"
$ws.Cells.Item(1622, 7).Value = "1,00"
$ws.Cells.Item(1622, 8).Value = "NUEVO"

$ws.Cells.Item(1623, 1).Value = "-"
$ws.Cells.Item(1623, 2).Value = "-"
$ws.Cells.Item(1623, 3).Value = "-"
$ws.Cells.Item(1623, 4).Value = "-"
$ws.Cells.Item(1623, 5).Value = "SUMINISTRO DE PIPExTE NIPPLE - LENGTH = 100 MM ASME B36.10 PExTE SEAMLESS ASTM A53 GR. B, 3/4 "", S-160"
$ws.Cells.Item(1623, 6).Value = "UN"
$ws.Cells.Item(1623, 7).Value = "3,00"
$ws.Cells.Item(1623, 8).Value = "NUEVO"

$ws.Cells.Item(1624, 1).Value = "-"
$ws.Cells.Item(1624, 2).Value = "-"
$ws.Cells.Item(1624, 3).Value = "-"
$ws.Cells.Item(1624, 4).Value = "-"
$ws.Cells.Item(1624, 5).Value = "SUMINISTRO DE REDUCER TEE, SCH 40 X SCH 30, BW, ASTM A-234 GR WPB, ASME B16.9 6""X 4"""
$ws.Cells.Item(1624, 6).Value = "UN"
$ws.Cells.Item(1624, 7).Value = "2,00"
$ws.Cells.Item(1624, 8).Value = "NUEVO"

$ws.Cells.Item(1625, 1).Value = "-"
$ws.Cells.Item(1625, 2).Value = "-"
$ws.Cells.Item(1625, 3).Value = "-"
$ws.Cells.Item(1625, 4).Value = "-"
$ws.Cells.Item(1625, 5).Value = "SUMINISTRO DE REDUCER TEE, SCH20 X SCH 30, BW, ASTM A-234 GR WPB, ASME B16.9 8""X 4"""
$ws.Cells.Item(1625, 6).Value = "UN"
$ws.Cells.Item(1625, 7).Value = "1,00"
$ws.Cells.Item(1625, 8).Value = "NUEVO"

$ws.Cells.Item(1626, 1).Value = "-"
$ws.Cells.Item(1626, 2).Value = "-"
$ws.Cells.Item(1626, 3).Value = "-"
$ws.Cells.Item(1626, 4).Value = "-"
$ws.Cells.Item(1626, 5).Value = "SUMINISTRO DE CONC. SWAGE, SCH 80 X SCH 80, PEXPE, ASTM A-234 GR WCB, MSS SP-95 1-1/2""X1"""
$ws.Cells.Item(1626, 6).Value = "UN"
$ws.Cells.Item(1626, 7).Value = "1,00"
$ws.Cells.Item(1626, 8).Value = "NUEVO"

$ws.Cells.Item(1627, 1).Value = "-"
$ws.Cells.Item(1627, 2).Value = "-"
$ws.Cells.Item(1627, 3).Value = "-"
$ws.Cells.Item(1627, 4).Value = "-"
$ws.Cells.Item(1627, 5).Value = "SUMINISTRO DE THREADOLET MSS SP-97 BW X SW 6000# FORGED ASTM A105, 6 "" X 1/2 """
$ws.Cells.Item(1627, 6).Value = "UN"
$ws.Cells.Item(1627, 7).Value = "1,00"
$ws.Cells.Item(1627, 8).Value = "NUEVO"

$ws.Cells.Item(1628, 1).Value = "-"
$ws.Cells.Item(1628, 2).Value = "-"
$ws.Cells.Item(1628, 3).Value = "-"
$ws.Cells.Item(1628, 4).Value = "-"
$ws.Cells.Item(1628, 5).Value = "SUMINISTRO DE WELDNECK FLANGE ASME B16.5 FLG X BW 150 LB FLAT FACE FORGED ASTM A105, 4 "", S-30"
$ws.Cells.Item(1628, 6).Value = "UN"
$ws.Cells.Item(1628, 7).Value = "4,00"
$ws.Cells.Item(1628, 8).Value = "NUEVO"

$ws.Cells.Item(1629, 1).Value = "-"
$ws.Cells.Item(1629, 2).Value = "-"
$ws.Cells.Item(1629, 3).Value = "-"
$ws.Cells.Item(1629, 4).Value = "-"
$ws.Cells.Item(1629, 5).Value = "SUMINISTRO DE WELDNECK FLANGE ASME B16.5 FLG X BW 150 LB FLAT FACE FORGED ASTM A105, 8 "", S-20"
$ws.Cells.Item(1629, 6).Value = "UN"
$ws.Cells.Item(1629, 7).Value = "4,00"
$ws.Cells.Item(1629, 8).Value = "NUEVO"

$ws.Cells.Item(1630, 1).Value = "-"
$ws.Cells.Item(1630, 2).Value = "-"
$ws.Cells.Item(1630, 3).Value = "-"
$ws.Cells.Item(1630, 4).Value = "-"
$ws.Cells.Item(1630, 5).Value = "SUMINISTRO DE WELDOLET MSS SP-97 BW X BW FORGED ASTM A105, 6"" X 2-1/2 "", S-40 X S-30"
$ws.Cells.Item(1630, 6).Value = "UN"
$ws.Cells.Item(1630, 7).Value = "4,00"
$ws.Cells.Item(1630, 8).Value = "NUEVO"

$ws.Cells.Item(1631, 1).Value = "-"
$ws.Cells.Item(1631, 2).Value = "-"
$ws.Cells.Item(1631, 3).Value = "-"
$ws.Cells.Item(1631, 4).Value = "-"
$ws.Cells.Item(1631, 5).Value = "PREFABRICACIÓN EN TALLER DE TUBERÍA DE ACERO AL CARBÓN DE Ø 2-1/2"" "
$ws.Cells.Item(1631, 6).Value = "-"
$ws.Cells.Item(1631, 7).Value = "15,28"
$ws.Cells.Item(1631, 8).Value = "NUEVO"

$ws.Cells.Item(1632, 1).Value = "-"
$ws.Cells.Item(1632, 2).Value = "-"
$ws.Cells.Item(1632, 3).Value = "-"
$ws.Cells.Item(1632, 4).Value = "-"
$ws.Cells.Item(1632, 5).Value = "PRUEBAS HIDROSTATICAS DE VÁLVULAS MANUALES BRIDADAS ANSI 150 DE DIAMETRO 2-1/2"""
$ws.Cells.Item(1632, 6).Value = "-"
$ws.Cells.Item(1632, 7).Value = "4,00"
$ws.Cells.Item(1632, 8).Value = "NUEVO"

$ws.Cells.Item(1633, 1).Value = "-"
$ws.Cells.Item(1633, 2).Value = "-"
$ws.Cells.Item(1633, 3).Value = "-"
$ws.Cells.Item(1633, 4).Value = "-"
$ws.Cells.Item(1633, 5).Value = "PRUEBAS HIDROSTATICAS DE VÁLVULAS MANUALES SW 800# DE DIAMETRO MENORES 2"""
$ws.Cells.Item(1633, 6).Value = "-"
$ws.Cells.Item(1633, 7).Value = "1,00"
$ws.Cells.Item(1633, 8).Value = "NUEVO"

$ws.Cells.Item(1634, 1).Value = "-"
$ws.Cells.Item(1634, 2).Value = "-"
$ws.Cells.Item(1634, 3).Value = "-"
$ws.Cells.Item(1634, 4).Value = "-"
$ws.Cells.Item(1634, 5).Value = "MONTAJE AÉREO DE TUBERÍA DE Ø 2-1/2”"
$ws.Cells.Item(1634, 6).Value = "-"
$ws.Cells.Item(1634, 7).Value = "5,00"
$ws.Cells.Item(1634, 8).Value = "NUEVO"

$ws.Cells.Item(1635, 1).Value = "-"
$ws.Cells.Item(1635, 2).Value = "-"
$ws.Cells.Item(1635, 3).Value = "-"
$ws.Cells.Item(1635, 4).Value = "-"
$ws.Cells.Item(1635, 5).Value = "MONTAJE DE VÁLVULAS Y CHEQUES DIAMETRO DE 2-1/2"""
$ws.Cells.Item(1635, 6).Value = "-"
$ws.Cells.Item(1635, 7).Value = "4,00"
$ws.Cells.Item(1635, 8).Value = "NUEVO"

